# Actualización automática: 2025-03-13 14:50:41
#
# Rows whose "Grupo" (column B) value moves from "LOTE1-General" to the new
# group "LOTE11-General". These are the exact rows touched by the edit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(47, 52, 56, 57, 58, 59, 60, 61, 62, 63, 64, 65)
foreach ($r in $rows) {
    $ws.Cells.Item($r, 2).Value = "LOTE11-General"
}

# Match the author's final on-screen view/selection state.
$ws.Range("B65").Select()
